$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows at row 22 (shifting the existing "Integrate with
# NodeServer and Express" ... "Post sales" block down by 2 rows, matching
# the blank-row-every-other-row layout used throughout the sheet).
$ws.Rows("22:23").Insert()

# Populate the newly inserted ToDo row (status left blank, same as the
# row layout convention: only B is filled, no C cell).
$ws.Range("B22").Value = "try to insert js lib references in the components rather than in index.html"

# Update the active selection to match the new target cell.
$ws.Range("B23").Select()
